$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148, shifting existing rows 148:175 down to 149:176
$ws.Rows("148").Insert()

# Populate the newly inserted row 148 with the new weekly price record
$ws.Range("A148").Value = 4
$ws.Range("B148").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C148").Value = "Los Lagos"
$ws.Range("D148").Value = 45209
$ws.Range("E148").Value = 10
$ws.Range("F148").Value = 100112031
$ws.Range("G148").Value = "Poroto verde"
$ws.Range("H148").Value = "Magnum"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 40
$ws.Range("K148").Value = 34000
$ws.Range("L148").Value = 34000
$ws.Range("M148").Value = 34000
$ws.Range("N148").Value = "$/malla 25 kilos"
$ws.Range("O148").Value = "Perú"
$ws.Range("P148").Value = 1360
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = "Hortaliza"
